# "code is optimized scroll function is pending"
# - Update the Pan No test data value (row 22 / column F) to a real-looking PAN.
# - Update the Terms & Conditions checkbox locator (row 25 / column D) from a
#   brittle indexed xpath to a className locator.
# - Leave the cursor positioned on the last edited cell (F22), scrolled so
#   column B is the first visible column, matching where the tester was
#   working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F22").Value = "AERPE9129F"
$ws.Range("D25").Value = "className = android.widget.CheckBox"

# Reposition the view: scroll so column B is left-most, then select F22
# (the cell that was last touched), matching the commit's view state.
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("F22").Select()
